$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.530.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.062.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.660'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.67'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -8.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("E10").Value = '  -7.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0748'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.108'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.904'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.358.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.122.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.448.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -12.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0856'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.18%  '
$ws.Range("E22").Value = '  -4.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  -5.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.35%  '
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0592'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.59%  '
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0806'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.87%  '
$ws.Range("E39").Value = '  -8.08%  '
$ws.Range("E40").Value = '  -5.76%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.06%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.42%  '
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0940'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.382.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -13.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.14%  '
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.246.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.19%  '
